# Add season record columns (Wins, Losses, Ties) to the NYM_2015 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - copy formatting from the existing header column (AC1) so the
# new headers share the same bold/centered/bordered style (s="1").
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record is the same for every player row (2 through 51): 90-72-0.
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 90
    $ws.Cells.Item($row, 31).Value = 72
    $ws.Cells.Item($row, 32).Value = 0
}
